{"js": "// The document starts with a single paragraph made of four runs:\n// \"Documentos \" + \"elaborados \" + \"en la \" + \"tarde\"\n// The edit collapses that paragraph down to one run of new text:\n// \"Subido commit en la rama masster a github\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n// Replacing the whole paragraph range rewrites its run(s) with a single\n// new run containing the given text (mirrors collapsing the 4 runs into 1).\nfirstParagraph.insertText(\"Subido commit en la rama masster a github\", \"Replace\");\nawait context.sync();\n", "ps1": "# The document starts with a single paragraph made of four runs:\n# \"Documentos \" + \"elaborados \" + \"en la \" + \"tarde\"\n# The edit collapses that paragraph down to one run of new text:\n# \"Subido commit en la rama masster a github\"\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n# wdReplaceAll = 2, wdFindContinue = 1\n$rng.Find.Execute(\n    \"Documentos elaborados en la tarde\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Subido commit en la rama masster a github\",\n    2\n)\n"}
